$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 636.8108
$ws.Range("I33").Value = 674.35297
$ws.Range("K33").Value = 674.35297
$ws.Range("M33").Value = -445.35297
$ws.Range("H38").Value = 2996.8
$ws.Range("I38").Value = 96
$ws.Range("J38").Value = 4240
$ws.Range("K38").Value = 288
$ws.Range("L38").Value = 12720
$ws.Range("M38").Value = 84
$ws.Range("N38").Value = -13464
$ws.Range("H39").Value = 171.1875
$ws.Range("I39").Value = 48.9
$ws.Range("J39").Value = 375
$ws.Range("K39").Value = 146.7
$ws.Range("L39").Value = 1125
$ws.Range("M39").Value = 149.3
$ws.Range("N39").Value = -1717
$ws.Range("H40").Value = 2073.087
$ws.Range("I40").Value = 1966.2727
$ws.Range("J40").Value = 2171
$ws.Range("K40").Value = 1966.2727
$ws.Range("L40").Value = 2171
$ws.Range("M40").Value = -1791.2727
$ws.Range("N40").Value = -2521
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H70").Value = 1676.3334
$ws.Range("I70").Value = 2300.25
$ws.Range("J70").Value = 963.2857
$ws.Range("K70").Value = 6900.75
$ws.Range("L70").Value = 2889.8571
$ws.Range("M70").Value = -6630.75
$ws.Range("N70").Value = -3429.8571
$ws.Range("H73").Value = 1676.3334
$ws.Range("I73").Value = 2300.25
$ws.Range("J73").Value = 963.2857
$ws.Range("K73").Value = 6900.75
$ws.Range("L73").Value = 2889.8571
$ws.Range("M73").Value = -5964.75
$ws.Range("N73").Value = -4761.8571
$ws.Range("H127").Value = 840.7143
$ws.Range("I127").Value = 552.1818
$ws.Range("J127").Value = 1158.1
$ws.Range("K127").Value = 1656.5454
$ws.Range("L127").Value = 3474.3
$ws.Range("M127").Value = 3303.4546
$ws.Range("N127").Value = -13394.3
$ws.Range("H132").Value = 5372.1724
$ws.Range("I132").Value = 5727.409
$ws.Range("J132").Value = 4255.7144
$ws.Range("K132").Value = 17182.227
$ws.Range("L132").Value = 12767.1432
$ws.Range("M132").Value = -14652.227
$ws.Range("N132").Value = -17827.1432
$ws.Range("H137").Value = 1304.2972
$ws.Range("I137").Value = 997.15
$ws.Range("J137").Value = 1665.6471
$ws.Range("K137").Value = 2991.45
$ws.Range("L137").Value = 4996.9413
$ws.Range("M137").Value = -441.4499999999998
$ws.Range("N137").Value = -10096.9413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 832.5833
$ws.Range("I2").Value = 839.0526
$ws.Range("J2").Value = 808
$ws.Range("K2").Value = 839.0526
$ws.Range("L2").Value = 808
$ws.Range("M2").Value = -726.0526
$ws.Range("N2").Value = -1034
$ws.Range("H32").Value = 16039.554
$ws.Range("I32").Value = 4916.8774
$ws.Range("J32").Value = 37840
$ws.Range("K32").Value = 4916.8774
$ws.Range("L32").Value = 37840
$ws.Range("M32").Value = -4629.8774
$ws.Range("N32").Value = -38414
$ws.Range("H45").Value = 2258.1538
$ws.Range("I45").Value = 2355.6
$ws.Range("J45").Value = 1933.3334
$ws.Range("K45").Value = 2355.6
$ws.Range("L45").Value = 1933.3334
$ws.Range("M45").Value = -1978.6
$ws.Range("N45").Value = -2687.3334
$ws.Range("H63").Value = 4716.6665
$ws.Range("I63").Value = 2466.6667
$ws.Range("K63").Value = 2466.6667
$ws.Range("M63").Value = -1780.6667
$ws.Range("H66").Value = 4716.6665
$ws.Range("I66").Value = 2466.6667
$ws.Range("K66").Value = 12333.3335
$ws.Range("M66").Value = -8901.333500000001
$ws.Range("H74").Value = 3453.5806
$ws.Range("I74").Value = 7835.3335
$ws.Range("J74").Value = 2401.96
$ws.Range("K74").Value = 7835.3335
$ws.Range("L74").Value = 2401.96
$ws.Range("M74").Value = -6961.3335
$ws.Range("N74").Value = -4149.96
$ws.Range("H77").Value = 3453.5806
$ws.Range("I77").Value = 7835.3335
$ws.Range("J77").Value = 2401.96
$ws.Range("K77").Value = 39176.6675
$ws.Range("L77").Value = 12009.8
$ws.Range("M77").Value = -34808.6675
$ws.Range("N77").Value = -20745.8
$ws.Range("H116").Value = 832.5833
$ws.Range("I116").Value = 839.0526
$ws.Range("J116").Value = 808
$ws.Range("K116").Value = 839.0526
$ws.Range("L116").Value = 808
$ws.Range("M116").Value = 1454.9474
$ws.Range("N116").Value = -5396
$ws.Range("H122").Value = 2112
$ws.Range("I122").Value = 1985.5385
$ws.Range("J122").Value = 2386
$ws.Range("K122").Value = 5956.6155
$ws.Range("L122").Value = 7158
$ws.Range("M122").Value = -3506.6155
$ws.Range("N122").Value = -12058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 832.5833
$ws.Range("I3").Value = 839.0526
$ws.Range("J3").Value = 808
$ws.Range("K3").Value = 839.0526
$ws.Range("L3").Value = 808
$ws.Range("M3").Value = -725.0526
$ws.Range("N3").Value = -1036
$ws.Range("H134").Value = 2664.138
$ws.Range("I134").Value = 2010.5333
$ws.Range("J134").Value = 3364.4285
$ws.Range("K134").Value = 6031.5999
$ws.Range("L134").Value = 10093.2855
$ws.Range("M134").Value = -3496.5999
$ws.Range("N134").Value = -15163.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2537.3
$ws.Range("I31").Value = 2032.3889
$ws.Range("J31").Value = 2821.3125
$ws.Range("K31").Value = 2032.3889
$ws.Range("L31").Value = 2821.3125
$ws.Range("M31").Value = -1737.3889
$ws.Range("N31").Value = -3411.3125
$ws.Range("H34").Value = 2537.3
$ws.Range("I34").Value = 2032.3889
$ws.Range("J34").Value = 2821.3125
$ws.Range("K34").Value = 2032.3889
$ws.Range("L34").Value = 2821.3125
$ws.Range("M34").Value = -1830.3889
$ws.Range("N34").Value = -3225.3125
$ws.Range("H58").Value = 1484.3077
$ws.Range("I58").Value = 921.06665
$ws.Range("J58").Value = 2252.3635
$ws.Range("K58").Value = 921.06665
$ws.Range("L58").Value = 2252.3635
$ws.Range("M58").Value = -718.06665
$ws.Range("N58").Value = -2658.3635
$ws.Range("H122").Value = 898.1053000000001
$ws.Range("I122").Value = 773.0769
$ws.Range("J122").Value = 1169
$ws.Range("K122").Value = 2319.2307
$ws.Range("L122").Value = 3507
$ws.Range("M122").Value = 130.7692999999999
$ws.Range("N122").Value = -8407
$ws.Range("H136").Value = 1484.3077
$ws.Range("I136").Value = 921.06665
$ws.Range("J136").Value = 2252.3635
$ws.Range("K136").Value = 2763.19995
$ws.Range("L136").Value = 6757.0905
$ws.Range("M136").Value = -213.1999500000002
$ws.Range("N136").Value = -11857.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4671
$ws.Range("J75").Value = 6000
$ws.Range("L75").Value = 18000
$ws.Range("N75").Value = -19996
$ws.Range("H78").Value = 4671
$ws.Range("J78").Value = 6000
$ws.Range("L78").Value = 54000
$ws.Range("N78").Value = -63984
$ws.Range("H107").Value = 270.3684
$ws.Range("I107").Value = 218.88889
$ws.Range("J107").Value = 316.7
$ws.Range("K107").Value = 656.6666700000001
$ws.Range("L107").Value = 950.0999999999999
$ws.Range("M107").Value = 1263.33333
$ws.Range("N107").Value = -4790.1
$ws.Range("H131").Value = 1667658.8
$ws.Range("I131").Value = 11111524
$ws.Range("J131").Value = 1094.2354
$ws.Range("K131").Value = 33334572
$ws.Range("L131").Value = 3282.7062
$ws.Range("M131").Value = -33329532
$ws.Range("N131").Value = -13362.7062
$ws.Range("H136").Value = 937.1667
$ws.Range("I136").Value = 726.9
$ws.Range("K136").Value = 2180.7
$ws.Range("M136").Value = 2919.3
$ws.Range("H140").Value = 1562.871
$ws.Range("I140").Value = 1176.2632
$ws.Range("J140").Value = 2175
$ws.Range("K140").Value = 3528.7896
$ws.Range("L140").Value = 6525
$ws.Range("M140").Value = 1651.2104
$ws.Range("N140").Value = -16885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 5274.6924
$ws.Range("J92").Value = 5274.6924
$ws.Range("L92").Value = 5274.6924
$ws.Range("N92").Value = -9018.6924
$ws.Range("H126").Value = 1652.1364
$ws.Range("I126").Value = 1338.0714
$ws.Range("J126").Value = 2201.75
$ws.Range("K126").Value = 4014.2142
$ws.Range("L126").Value = 6605.25
$ws.Range("M126").Value = -1544.2142
$ws.Range("N126").Value = -11545.25
$ws.Range("H132").Value = 2240.8948
$ws.Range("I132").Value = 1480
$ws.Range("K132").Value = 4440
$ws.Range("M132").Value = -1910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1355.2667
$ws.Range("I7").Value = 735.36365
$ws.Range("J7").Value = 3060
$ws.Range("K7").Value = 735.36365
$ws.Range("L7").Value = 3060
$ws.Range("M7").Value = -623.36365
$ws.Range("N7").Value = -3284
$ws.Range("H68").Value = 240350.23
$ws.Range("I68").Value = 834348.3
$ws.Range("J68").Value = 2751
$ws.Range("K68").Value = 834348.3
$ws.Range("L68").Value = 2751
$ws.Range("M68").Value = -833599.3
$ws.Range("N68").Value = -4249
$ws.Range("H71").Value = 240350.23
$ws.Range("I71").Value = 834348.3
$ws.Range("J71").Value = 2751
$ws.Range("K71").Value = 4171741.5
$ws.Range("L71").Value = 13755
$ws.Range("M71").Value = -4167997.5
$ws.Range("N71").Value = -21243
$ws.Range("H122").Value = 3373
$ws.Range("I122").Value = 3671.2856
$ws.Range("J122").Value = 2851
$ws.Range("K122").Value = 11013.8568
$ws.Range("L122").Value = 8553
$ws.Range("M122").Value = -8563.856800000001
$ws.Range("N122").Value = -13453
$ws.Range("H126").Value = 1355.2667
$ws.Range("I126").Value = 735.36365
$ws.Range("J126").Value = 3060
$ws.Range("K126").Value = 2206.09095
$ws.Range("L126").Value = 9180
$ws.Range("M126").Value = 263.9090500000002
$ws.Range("N126").Value = -14120
$ws.Range("H136").Value = 3932.889
$ws.Range("I136").Value = 5490.3706
$ws.Range("J136").Value = 1596.6666
$ws.Range("K136").Value = 16471.1118
$ws.Range("L136").Value = 4789.9998
$ws.Range("M136").Value = -13921.1118
$ws.Range("N136").Value = -9889.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1050
$ws.Range("I96").Value = 1150
$ws.Range("J96").Value = 850
$ws.Range("K96").Value = 1150
$ws.Range("L96").Value = 850
$ws.Range("M96").Value = 223
$ws.Range("N96").Value = -3596
$ws.Range("H122").Value = 1018.7647
$ws.Range("I122").Value = 827.0526
$ws.Range("J122").Value = 1261.6
$ws.Range("K122").Value = 2481.1578
$ws.Range("L122").Value = 3784.8
$ws.Range("M122").Value = -31.15779999999995
$ws.Range("N122").Value = -8684.799999999999
$ws.Range("H126").Value = 5311.5557
$ws.Range("I126").Value = 6400.5713
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 19201.7139
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -16731.7139
$ws.Range("N126").Value = -9440
